$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 5.8
$ws.Range("C2").Value = 12.9
$ws.Range("B3").Value = 5.7
$ws.Range("C3").Value = 11.7
$ws.Range("C5").Value = 29
